$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new trailing columns to the export/import header row
$ws.Range("CR1").Value = "Jurisdiction Path"
$ws.Range("CS1").Value = "Group Number"

# Match the formatting used for the other header cells (thin border, centered,
# wrapped text, no fill) - same look as the rest of row 1's header cells.
$rng = $ws.Range("CR1:CS1")
$rng.HorizontalAlignment = -4108  # xlCenter
$rng.VerticalAlignment = -4108    # xlCenter
$rng.WrapText = $true
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# Leave selection on the next empty cell, as in the authored workbook
$ws.Range("CS2").Select()
